$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks so they can be rebuilt in the correct order
$ws.Hyperlinks.Delete()

# Fill in the repeating credential rows 4-12 (pattern matches rows 2/3)
for ($r = 4; $r -le 12; $r++) {
  if ($r % 2 -eq 0) {
    $ws.Range("A$r").Value = "dineshkumar.icon@gmail.com"
  } else {
    $ws.Range("A$r").Value = "dineshkumar.icon.dk@gmail.com"
  }
  $ws.Range("B$r").Value = "Dinnu@247"
}

# Recreate hyperlinks for the original cells first (rId1-rId3)
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Dinnu@247") | Out-Null
$ws.Range("B2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
$ws.Range("A3").Style = "Hyperlink"

# Column A hyperlinks for the new rows: even rows (gmail.com) first, then odd rows (icon.dk)
foreach ($r in 4,6,8,10,12) {
  $ws.Hyperlinks.Add($ws.Range("A$r"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
  $ws.Range("A$r").Style = "Hyperlink"
}
foreach ($r in 5,7,9,11) {
  $ws.Hyperlinks.Add($ws.Range("A$r"), "mailto:dineshkumar.icon.dk@gmail.com") | Out-Null
  $ws.Range("A$r").Style = "Hyperlink"
}

# A single hyperlink spanning B3:B12, then individual ones for each cell in that range
$ws.Hyperlinks.Add($ws.Range("B3:B12"), "mailto:Dinnu@247", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Dinnu@247") | Out-Null
foreach ($r in 3,4,5,6,7,8,9,10,11,12) {
  $ws.Hyperlinks.Add($ws.Range("B$r"), "mailto:Dinnu@247") | Out-Null
  $ws.Range("B$r").Style = "Hyperlink"
}

# Update selection to match target
$ws.Range("B12").Select()

Write-Host "Final hyperlink count:" $ws.Hyperlinks.Count
